# Scheduled runner update: refresh market-board derived columns (H-N)
# across the Leve profitability sheets with newly scraped prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 385.7143
$ws.Range("J12").Value = 128
$ws.Range("L12").Value = 128
$ws.Range("N12").Value = -468
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("H98").Value = 1431.25
$ws.Range("I98").Value = 1453.9474
$ws.Range("K98").Value = 1453.9474
$ws.Range("M98").Value = 44.05259999999998
$ws.Range("H122").Value = 1431.25
$ws.Range("I122").Value = 1453.9474
$ws.Range("K122").Value = 4361.8422
$ws.Range("M122").Value = -1911.8422
$ws.Range("H129").Value = 898.6892
$ws.Range("I129").Value = 1219.8
$ws.Range("J129").Value = 875.4203
$ws.Range("K129").Value = 3659.4
$ws.Range("L129").Value = 2626.2609
$ws.Range("M129").Value = 1340.6
$ws.Range("N129").Value = -12626.2609
$ws.Range("H137").Value = 1552.6154
$ws.Range("J137").Value = 2500
$ws.Range("L137").Value = 7500
$ws.Range("N137").Value = -12600
$ws.Range("H140").Value = 53083.41
$ws.Range("J140").Value = 53083.41
$ws.Range("L140").Value = 53083.41
$ws.Range("N140").Value = -63443.41
$ws.Range("H141").Value = 1753839
$ws.Range("I141").Value = 2802733.8
$ws.Range("K141").Value = 8408201.399999999
$ws.Range("M141").Value = -8403021.399999999
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3662.3044
$ws.Range("I32").Value = 2390.95
$ws.Range("K32").Value = 2390.95
$ws.Range("M32").Value = -2103.95
$ws.Range("H55").Value = 14280
$ws.Range("J55").Value = 14280
$ws.Range("L55").Value = 14280
$ws.Range("N55").Value = -14910
$ws.Range("H61").Value = 3153.8823
$ws.Range("I61").Value = 2551.4333
$ws.Range("J61").Value = 7672.25
$ws.Range("K61").Value = 2551.4333
$ws.Range("L61").Value = 7672.25
$ws.Range("M61").Value = -2339.4333
$ws.Range("N61").Value = -8096.25
$ws.Range("H122").Value = 3134.5557
$ws.Range("I122").Value = 3134.5557
$ws.Range("K122").Value = 9403.667099999999
$ws.Range("M122").Value = -6953.667099999999
$ws.Range("H132").Value = 1606.9062
$ws.Range("I132").Value = 1045.6428
$ws.Range("J132").Value = 2043.4445
$ws.Range("K132").Value = 3136.9284
$ws.Range("L132").Value = 6130.333500000001
$ws.Range("M132").Value = -606.9284000000002
$ws.Range("N132").Value = -11190.3335
$ws.Range("H136").Value = 3153.8823
$ws.Range("I136").Value = 2551.4333
$ws.Range("J136").Value = 7672.25
$ws.Range("K136").Value = 7654.2999
$ws.Range("L136").Value = 23016.75
$ws.Range("M136").Value = -5104.2999
$ws.Range("N136").Value = -28116.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1375.5294
$ws.Range("I20").Value = 1370.6522
$ws.Range("J20").Value = 1385.7273
$ws.Range("K20").Value = 1370.6522
$ws.Range("L20").Value = 1385.7273
$ws.Range("M20").Value = -1123.6522
$ws.Range("N20").Value = -1879.7273
$ws.Range("H86").Value = 75662.44500000001
$ws.Range("I86").Value = 1402.1904
$ws.Range("J86").Value = 335573.34
$ws.Range("K86").Value = 1402.1904
$ws.Range("L86").Value = 335573.34
$ws.Range("M86").Value = -279.1904
$ws.Range("N86").Value = -337819.34
$ws.Range("H89").Value = 75662.44500000001
$ws.Range("I89").Value = 1402.1904
$ws.Range("J89").Value = 335573.34
$ws.Range("K89").Value = 7010.951999999999
$ws.Range("L89").Value = 1677866.7
$ws.Range("M89").Value = -1394.951999999999
$ws.Range("N89").Value = -1689098.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 246.45454
$ws.Range("I7").Value = 154
$ws.Range("K7").Value = 154
$ws.Range("M7").Value = -41
$ws.Range("H31").Value = 2796.4473
$ws.Range("I31").Value = 1837.25
$ws.Range("K31").Value = 1837.25
$ws.Range("M31").Value = -1542.25
$ws.Range("H34").Value = 2796.4473
$ws.Range("I34").Value = 1837.25
$ws.Range("K34").Value = 1837.25
$ws.Range("M34").Value = -1635.25
$ws.Range("H99").Value = 771793.3
$ws.Range("I99").Value = 2501424.8
$ws.Range("J99").Value = 3068.2222
$ws.Range("K99").Value = 2501424.8
$ws.Range("L99").Value = 3068.2222
$ws.Range("M99").Value = -2499926.8
$ws.Range("N99").Value = -6064.2222
$ws.Range("H105").Value = 1317.375
$ws.Range("I105").Value = 1218.2858
$ws.Range("K105").Value = 1218.2858
$ws.Range("M105").Value = 528.7141999999999
$ws.Range("H122").Value = 2718.7693
$ws.Range("I122").Value = 2830
$ws.Range("J122").Value = 2623.4285
$ws.Range("K122").Value = 8490
$ws.Range("L122").Value = 7870.2855
$ws.Range("M122").Value = -6040
$ws.Range("N122").Value = -12770.2855
$ws.Range("H126").Value = 771793.3
$ws.Range("I126").Value = 2501424.8
$ws.Range("J126").Value = 3068.2222
$ws.Range("K126").Value = 7504274.399999999
$ws.Range("L126").Value = 9204.6666
$ws.Range("M126").Value = -7501804.399999999
$ws.Range("N126").Value = -14144.6666
$ws.Range("H132").Value = 2266.2307
$ws.Range("I132").Value = 1344.6428
$ws.Range("J132").Value = 3341.4167
$ws.Range("K132").Value = 4033.9284
$ws.Range("L132").Value = 10024.2501
$ws.Range("M132").Value = -1503.9284
$ws.Range("N132").Value = -15084.2501
$ws.Range("H134").Value = 952.86664
$ws.Range("I134").Value = 941.1667
$ws.Range("J134").Value = 999.6667
$ws.Range("K134").Value = 2823.5001
$ws.Range("L134").Value = 2999.0001
$ws.Range("M134").Value = -288.5001000000002
$ws.Range("N134").Value = -8069.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 809.2222
$ws.Range("I5").Value = 751
$ws.Range("J5").Value = 838.3333
$ws.Range("K5").Value = 2253
$ws.Range("L5").Value = 2514.9999
$ws.Range("M5").Value = -2141
$ws.Range("N5").Value = -2738.9999
$ws.Range("H60").Value = 1941.6666
$ws.Range("I60").Value = 1941.6666
$ws.Range("K60").Value = 5824.9998
$ws.Range("M60").Value = -5573.9998
$ws.Range("H81").Value = 28910444
$ws.Range("I81").Value = 996.6667
$ws.Range("J81").Value = 39751490
$ws.Range("K81").Value = 2990.0001
$ws.Range("L81").Value = 119254470
$ws.Range("M81").Value = -1867.0001
$ws.Range("N81").Value = -119256716
$ws.Range("H84").Value = 28910444
$ws.Range("I84").Value = 996.6667
$ws.Range("J84").Value = 39751490
$ws.Range("K84").Value = 8970.0003
$ws.Range("L84").Value = 357763410
$ws.Range("M84").Value = -3354.0003
$ws.Range("N84").Value = -357774642
$ws.Range("H107").Value = 906.3333
$ws.Range("J107").Value = 906.3333
$ws.Range("L107").Value = 2718.9999
$ws.Range("N107").Value = -6558.9999
$ws.Range("H130").Value = 2100
$ws.Range("J130").Value = 2100
$ws.Range("L130").Value = 6300
$ws.Range("N130").Value = -16340
$ws.Range("H131").Value = 9627.253000000001
$ws.Range("J131").Value = 10642.226
$ws.Range("L131").Value = 31926.678
$ws.Range("N131").Value = -42006.678
$ws.Range("H135").Value = 809.2222
$ws.Range("I135").Value = 751
$ws.Range("J135").Value = 838.3333
$ws.Range("K135").Value = 6759
$ws.Range("L135").Value = 7544.9997
$ws.Range("M135").Value = -4224
$ws.Range("N135").Value = -12614.9997
$ws.Range("H137").Value = 3236.3809
$ws.Range("I137").Value = 1971.7142
$ws.Range("J137").Value = 3868.7144
$ws.Range("K137").Value = 5915.142599999999
$ws.Range("L137").Value = 11606.1432
$ws.Range("M137").Value = -815.1425999999992
$ws.Range("N137").Value = -21806.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 19920
$ws.Range("J46").Value = 19920
$ws.Range("L46").Value = 19920
$ws.Range("N46").Value = -20232
$ws.Range("H122").Value = 2131.7144
$ws.Range("J122").Value = 2301.6
$ws.Range("L122").Value = 6904.799999999999
$ws.Range("N122").Value = -11804.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6831
$ws.Range("J40").Value = 11674.5
$ws.Range("L40").Value = 11674.5
$ws.Range("N40").Value = -11946.5
$ws.Range("H46").Value = 2249.75
$ws.Range("I46").Value = 1133.3334
$ws.Range("K46").Value = 1133.3334
$ws.Range("M46").Value = -945.3334
$ws.Range("H133").Value = 79000
$ws.Range("J133").Value = 79000
$ws.Range("L133").Value = 79000
$ws.Range("N133").Value = -84060
$ws.Range("H136").Value = 2551.5
$ws.Range("I136").Value = 1217.5385
$ws.Range("J136").Value = 4478.3335
$ws.Range("K136").Value = 3652.6155
$ws.Range("L136").Value = 13435.0005
$ws.Range("M136").Value = -1102.6155
$ws.Range("N136").Value = -18535.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 53285.285
$ws.Range("J108").Value = 53285.285
$ws.Range("L108").Value = 53285.285
$ws.Range("N108").Value = -60965.285
$ws.Range("H122").Value = 638888.5
$ws.Range("I122").Value = 638888.5
$ws.Range("K122").Value = 1916665.5
$ws.Range("M122").Value = -1914215.5
$ws.Range("H126").Value = 20055.777
$ws.Range("I126").Value = 32625.75
$ws.Range("J126").Value = 9999.799999999999
$ws.Range("K126").Value = 97877.25
$ws.Range("L126").Value = 29999.4
$ws.Range("M126").Value = -95407.25
$ws.Range("N126").Value = -34939.39999999999
$ws.Range("H136").Value = 17923446
$ws.Range("I136").Value = 25254624
$ws.Range("J136").Value = 2783.2222
$ws.Range("K136").Value = 75763872
$ws.Range("L136").Value = 8349.6666
$ws.Range("M136").Value = -75761322
$ws.Range("N136").Value = -13449.6666
